$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.121.69'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '2.940.54'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '376.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.53%  '
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -2.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.78'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = '3.401.00'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.23%  '
$ws.Range("E15").Value = '  -2.36%  '
$ws.Range("D16").Value = '2.940.90'
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '51.065.19'
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("E19").Value = '  -6.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.16%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '263.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").Value = '  +2.48%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.87%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.24%  '
$ws.Range("E28").Value = '  +5.08%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("E32").Value = '  -0.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '50.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("E36").Value = '  -3.96%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  -5.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.32%  '
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("E42").Value = '  -4.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '121.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.272'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.92%  '
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").Value = '2.002.32'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.41%  '
